$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bug fix: Material values were swapped between row 4 (CageID 15A) and row 8 (CageID 32) ---
$ws.Range("E4").Value = "Metal"
$ws.Range("E8").Value = "Wood"

# --- Normalize formatting (center alignment) on rows 35-39 which were missing the style ---
$fmtRange = $ws.Range("A35:E39")
$fmtRange.HorizontalAlignment = -4108
$fmtRange.VerticalAlignment = -4108

# --- Add new cage records (rows 40-42) ---
$ws.Range("A40").Value = 43
$ws.Range("B40").Value = 0.5
$ws.Range("C40").Value = 2
$ws.Range("D40").Value = 5.222
$ws.Range("E40").Value = "Plastic"

$ws.Range("A41").Value = "423FF"
$ws.Range("B41").Value = 3.2
$ws.Range("C41").Value = 33
$ws.Range("D41").Value = 3
$ws.Range("E41").Value = "Plastic"

$ws.Range("A42").Value = "187A"
$ws.Range("B42").Value = 15.8
$ws.Range("C42").Value = 10
$ws.Range("D42").Value = 30
$ws.Range("E42").Value = "Wood"

$newRange = $ws.Range("A40:E42")
$newRange.HorizontalAlignment = -4108
$newRange.VerticalAlignment = -4108

# --- Column widths / default formatting for the sheet ---
$ws.Columns.Item(1).ColumnWidth = 11.7

# --- Update selection to match where the user ended up ---
$ws.Range("G41").Select()
